# Commit: "Updates to code and figures"
#
# The author removed the second slide (sldId="257", the classification
# scenario table/figure slide) from FigureS1.pptx, keeping only the first
# slide (sldId="256"). Everything else in the diff (the refreshed
# "datetimeFigureOut" date-placeholder caches on the slide layouts/master,
# and the Office change-tracking metadata in ppt/changesInfos) is just the
# automatic side effect of PowerPoint re-saving the deck on a later date --
# not a deliberate, scriptable content edit.

$p = $ppt.ActivePresentation

# Slide 2 is the one with sldId="257" (rId3 -> ppt/slides/slide2.xml) --
# delete it, leaving only slide 1 (sldId="256").
$p.Slides.Item(2).Delete()
